$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Work from the bottom of the sheet upward so that row numbers for the
# rows we haven't touched yet stay stable.

# 1) Remove the old TATIANA row (005348011 / 80.91) that sat right after
#    003115072 VICTOR (row 70).
$ws.Rows.Item(70).Delete()

# 2) Remove the old HEPTA row (004359408 / 86.35) that sat right before
#    004584982 BRUNO (row 64).
$ws.Rows.Item(64).Delete()

# 3) Insert a new TATIANA row (005348011 / 4048.51) right before
#    004467884 ANA (row 13).
$ws.Rows.Item(13).Insert()
$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "005348011"
$ws.Cells.Item(13, 2).Value = "TATIANA"
$ws.Cells.Item(13, 3).Value = 4048.51

# 4) Insert a new ADELE row (004575632 / 11000) right before
#    004748761 MARCELO (row 12).
$ws.Rows.Item(12).Insert()
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "004575632"
$ws.Cells.Item(12, 2).Value = "ADELE"
$ws.Cells.Item(12, 3).Value = 11000

# 5) Remove the old CLOTILDE row (005639338 / 20000) that sat right after
#    005064129 THIAGO (row 8).
$ws.Rows.Item(8).Delete()

# 6) Insert the new HEPTA (004359408 / 83216.35) and CLOTILDE
#    (005639338 / 30000) rows right before 004224011 THOMAS (row 6).
$ws.Rows.Item(6).Insert()
$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "004359408"
$ws.Cells.Item(6, 2).Value = "HEPTA"
$ws.Cells.Item(6, 3).Value = 83216.35

$ws.Rows.Item(7).Insert()
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "005639338"
$ws.Cells.Item(7, 2).Value = "CLOTILDE"
$ws.Cells.Item(7, 3).Value = 30000
